$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before the current row 601, shifting the existing
# rows 601-609 down to 604-612 (dimension grows from R609 to R612).
$ws.Rows("601:603").Insert()

# New weekly price data (row, then columns A..R)
$rowsData = @(
    @(3, "Femacal de La Calera", "Coquimbo", 44939, 5, 100112028, "Sandia", "Sin especificar", "Extra",   450, 3000, 3000, 3000, "`$/unidad", "Región Metropolitana", 3000, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44939, 5, 100112028, "Sandia", "Sin especificar", "Primera", 710, 2000, 2300, 2097, "`$/unidad", "Región Metropolitana", 2097, 1, "Hortaliza"),
    @(3, "Femacal de La Calera", "Coquimbo", 44939, 5, 100112028, "Sandia", "Sin especificar", "Segunda", 700, 1500, 1700, 1571, "`$/unidad", "Región Metropolitana", 1571, 1, "Hortaliza")
)

$startRow = 601
for ($i = 0; $i -lt $rowsData.Count; $i++) {
    $row = $rowsData[$i]
    $rowNum = $startRow + $i
    for ($j = 0; $j -lt $row.Count; $j++) {
        $colNum = $j + 1
        $ws.Cells.Item($rowNum, $colNum).Value = $row[$j]
    }
}
